$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "95.362.55"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.66%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.625.87"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -2.30%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "2.35"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +24.58%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.00"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.28%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "226.54"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -4.55%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "640.02"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -1.94%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.419"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.35%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.12"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +5.40%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.00"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.14%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.623.35"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -2.26%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "47.63"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +7.66%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.208"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.66%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000295"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -2.85%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.54"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.31%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.306.07"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.18%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "95.192.91"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -1.56%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.82"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.72%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.636.37"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.62%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "19.48"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +3.97%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.82"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.55%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.522"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.27%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "517.41"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.99%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.26"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -4.56%  "

$ws.Range("B25").Value = "Hedera"
$ws.Range("C25").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.243"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +26.54%  "

$ws.Range("B26").Value = "Litecoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "120.09"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +18.43%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000205"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -3.03%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.77"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.40%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "12.78"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -4.52%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "12.80"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +5.32%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.93"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -2.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.00"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +0.06%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.181"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -4.39%  "

$ws.Range("B34").Value = "Binance-PegBSC-USD"
$ws.Range("C34").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +0.08%  "

$ws.Range("B35").Value = "Fetch.AI"
$ws.Range("C35").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.78"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.99%  "

$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "32.11"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.69%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.592"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -1.25%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "602.52"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -7.52%  "

$ws.Range("E39").Value = "  -0.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.37"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -5.35%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.91"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +1.88%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.488"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +9.56%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.159"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -0.88%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.79"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.81%  "

$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0484"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +6.52%  "

$ws.Range("B46").Value = "ImmutableX"
$ws.Range("C46").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.94"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -4.33%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.931"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.19%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.77"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.00%  "

$ws.Range("B49").Value = "WhiteBITCoin"
$ws.Range("C49").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.46"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.62%  "

$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.23"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -2.49%  "

$ws.Range("B51").Value = "OKB"
$ws.Range("C51").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "54.28"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.76%  "
